$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, copying the formatting used by the rest of the
# header row (B1:G1) so it picks up the same style index.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value in H2, below the new header.
$ws.Range("H2").Value = 1
